$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scalability subSeq")
$ws.Range("A1").Value = "TEST"
